# DGS time-series update: append the 2021/08/23 report row (row 71).
#
# Column A stores the report date as literal text (it is a shared string
# in the workbook even though the column's display format is yyyy/mm/dd),
# so we briefly force a Text number format while assigning the value to
# stop Excel from auto-converting the "2021/08/23" literal into a date
# serial number, then restore the column's normal date display format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "2021/08/23"
$ws.Range("A71").NumberFormat = "yyyy/mm/dd"

$ws.Range("B71").Value = 310.4
$ws.Range("C71").Value = 314.6
$ws.Range("D71").Value = 0.98
$ws.Range("E71").Value = 0.98

# Move the selection to A72, matching where Excel leaves the cursor after
# entering data in A71 (mirrors the post-edit active cell in the workbook).
$ws.Range("A72").Select() | Out-Null
